$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number format on the whole target range so numeric-looking
# strings (e.g. '67', '8168') and empty strings are preserved as text
# cells instead of being coerced to numbers or dropped entirely.
$ws.Range("A1:G3").NumberFormat = "@"

$ws.Range("A1").Value = 'Amazfit Active Akıllı Saat (Amazfit Türkiye Garantili)'
$ws.Range("B1").Value = ""
$ws.Range("C1").Value = ""
$ws.Range("D1").Value = 'https://n11scdn.akamaized.net/a1/602_857/09/67/05/91/IMG-1977890781065674594.jpg'
$ws.Range("E1").Value = '[''https://n11scdn.akamaized.net/a1/602_857/09/67/05/91/IMG-1977890781065674594.jpg'', ''https://n11scdn.akamaized.net/a1/60_86/09/67/05/91/IMG-1977890781065674594.jpg'', ''https://n11scdn.akamaized.net/a1/60_86/15/25/87/43/IMG-5803356974882669716.jpg'', ''https://n11scdn.akamaized.net/a1/60_86/09/43/32/64/IMG-7912115531430129761.jpg'', ''https://n11scdn.akamaized.net/a1/60_86/13/08/46/93/IMG-887897739486994108.jpg'', ''https://n11scdn.akamaized.net/a1/60_86/15/37/66/19/IMG-2094926088427494656.jpg'', ''https://n11scdn.akamaized.net/a1/60_86/03/40/38/24/IMG-422236255780343523.jpg'']'
$ws.Range("F1").Value = ""
$ws.Range("G1").Value = ""

$ws.Range("A2").Value = 'Samsung Galaxy Tab A9+ Wi-Fi SM-X210 8 GB 128 GB 11" Tablet'
$ws.Range("B2").Value = ""
$ws.Range("C2").Value = ""
$ws.Range("D2").Value = 'https://n11scdn.akamaized.net/a1/602_857/10/00/94/73/IMG-3886570089543208949.jpg'
$ws.Range("E2").Value = '[''https://n11scdn.akamaized.net/a1/602_857/10/00/94/73/IMG-3886570089543208949.jpg'', ''https://n11scdn.akamaized.net/a1/60_86/10/00/94/73/IMG-3886570089543208949.jpg'', ''https://n11scdn.akamaized.net/a1/60_86/07/30/60/21/IMG-8525196323097482019.jpg'', ''https://n11scdn.akamaized.net/a1/60_86/16/32/70/81/IMG-2575484217061974956.jpg'', ''https://n11scdn.akamaized.net/a1/60_86/09/00/03/42/IMG-7057197388026173745.jpg'', ''https://n11scdn.akamaized.net/a1/60_86/06/45/35/09/IMG-5919884287401727627.jpg'', ''https://n11scdn.akamaized.net/a1/60_86/14/00/73/40/IMG-8217500455611729157.jpg'', ''https://n11scdn.akamaized.net/a1/500_700/24/02/16/74/20/86/91/76/06/26/05/54/86651391626104357813.png'', ''https://n11scdn.akamaized.net/a1/500_700/24/02/16/62/80/80/79/32/08/78/38/14/92412933804080057162.png'', ''https://n11scdn.akamaized.net/a1/500_700/24/02/16/81/07/39/21/28/39/63/27/84/71086975294612718064.png'', ''https://n11scdn.akamaized.net/a1/500_700/24/02/09/48/31/90/25/04/04/21/94/49/31917057263205894343.png'']'
$ws.Range("F2").Value = '4,5'
$ws.Range("G2").Value = '67'

$ws.Range("A3").Value = 'Xiaomi Mi Wifi Pro 300 Mbps 2.4 GHz Sinyal Güçlendirici (Türkiye Versiyon)'
$ws.Range("B3").Value = ""
$ws.Range("C3").Value = ""
$ws.Range("D3").Value = 'https://n11scdn.akamaized.net/a1/602_857/05/53/69/25/IMG-8811551717331288430.jpg'
$ws.Range("E3").Value = '[''https://n11scdn.akamaized.net/a1/602_857/05/53/69/25/IMG-8811551717331288430.jpg'', ''https://n11scdn.akamaized.net/a1/60_86/05/53/69/25/IMG-8811551717331288430.jpg'', ''https://n11scdn.akamaized.net/a1/60_86/07/37/79/79/IMG-2780517974147884701.jpg'', ''https://n11scdn.akamaized.net/a1/60_86/02/00/95/36/82946241.jpg'', ''https://n11scdn.akamaized.net/a1/60_86/10/67/88/52/IMG-4204343210519024103.jpg'', ''https://n11scdn.akamaized.net/a1/60_86/15/31/87/26/IMG-5272898711302070857.jpg'', ''https://n11scdn.akamaized.net/a1/500_700/24/02/23/58/97/55/05/84/44/66/61/08/95371355153595515091.jpg'', ''https://n11scdn.akamaized.net/a1/500_700/24/02/23/91/42/47/78/15/63/37/13/49/76186264513104557719.jpg'', ''https://n11scdn.akamaized.net/a1/500_700/24/02/22/78/47/90/19/93/38/56/20/71/23358246157126711727.png'', ''https://n11scdn.akamaized.net/a1/500_700/24/02/22/52/97/97/95/08/50/87/08/46/25891089956506912282.png'', ''https://n11scdn.akamaized.net/a1/500_700/24/02/22/92/06/22/49/18/43/74/50/04/36802597043987331134.png'']'
$ws.Range("F3").Value = '4,5'
$ws.Range("G3").Value = '8168'

